$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: tag existing "Player"-side entry with category "Sprites" (N11) ---
$ws.Range("N11").Value = "Sprites"

# --- Row 12: same tag added to N12 (new cell, needs the bordered "N"-column style) ---
$ws.Range("N8").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = "Sprites"

# --- Row 13: fill in the Viki-side entry (J/K/L/M) that was previously blank, plus tag N13 ---
$ws.Range("J13").Value = 45363
$ws.Range("K13").Formula = "=8"
$ws.Range("L13").Formula = "=9+35/60"

$ws.Range("N8").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = "Sprites"

# --- Rows 14-27: add the blank, date-formatted J cell (matches the rest of the J column) ---
$ws.Range("J12").Copy()
$jRows = 14..27
foreach ($r in $jRows) {
  $ws.Range("J$r").PasteSpecial(-4122)
}

# --- New time-log entry, row 31 ---
$ws.Range("A26").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = 45363

$ws.Range("B31").Formula = "=8"
$ws.Range("C31").Formula = "=9+35/60"

$ws.Range("D26").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Formula = "=C31-B31"

$ws.Range("E27").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = "Main menu, jumppad and spike"

# --- selection cursor moved to the next empty row, as left by the editing user ---
$ws.Range("E32").Select()

Write-Output "edit applied"
